$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New purchase record: AGC, bought 3 Mar 2020
$ws.Range("A5").Value = "AGC"
$ws.Range("B5").Value = 390.95
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "3 Mar 2020"
$ws.Range("D5").Value = 8

# New purchase record: DA01, bought 3 Mar 2020
$ws.Range("A6").Value = "DA01"
$ws.Range("B6").Value = 458.6
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "3 Mar 2020"
$ws.Range("D6").Value = 7

# Page setup used when the sheet was re-saved
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on the first empty row below the table
$ws.Range("A7").Select() | Out-Null
